# Update inventory updated with new sales script runs for early 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warehouse")

# Row 2: ccao-condominium-pin_condo_char.R -> ran 2/22/2023, clear the "waiting on valuations" note
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = #2/22/2023#
$ws.Range("A2:C2").Style = "Normal"

# Row 13: spatial-ccao-county.R -> no longer highlighted (plain formatting)
$ws.Range("A13:C13").Style = "Normal"
$ws.Range("B13").Value = #1/12/2022#
$ws.Range("B13").NumberFormat = "m/d/yyyy"

# Row 16: ccao-condominium_parking.R -> re-run 2/27/2023, clear the "needs to be run again" note
$ws.Range("B16").Value = #2/27/2023#
$ws.Range("D16").ClearContents()

# Row 17: ccao-dictionary.R -> re-run 2/27/2023, clear the "needs to be run again" note
$ws.Range("B17").Value = #2/27/2023#
$ws.Range("D17").ClearContents()

# Move selection to B17 to match the saved cursor position
$ws.Range("B17").Select()
